$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 351, pushing existing rows 351-404 down to 352-405
$ws.Rows("351:351").Insert()

# Populate the newly inserted row 351 with the new data entry
$ws.Cells.Item(351, 1).Value = 5
$ws.Cells.Item(351, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(351, 3).Value = "Maule"
$ws.Cells.Item(351, 4).Value = 44505
$ws.Cells.Item(351, 5).Value = 7
$ws.Cells.Item(351, 6).Value = 100112004
$ws.Cells.Item(351, 7).Value = "Cebolla"
$ws.Cells.Item(351, 8).Value = "Sin especificar"
$ws.Cells.Item(351, 9).Value = "1a nueva(o)"
$ws.Cells.Item(351, 10).Value = 60000
$ws.Cells.Item(351, 11).Value = 1200
$ws.Cells.Item(351, 12).Value = 1200
$ws.Cells.Item(351, 13).Value = 1200
$ws.Cells.Item(351, 14).Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Cells.Item(351, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(351, 16).Value = 120
$ws.Cells.Item(351, 17).Value = 10
$ws.Cells.Item(351, 18).Value = "Hortaliza"
